# Auto update stock data
# Update the "Date_1" column (A) value "2025/11/30" -> "2025/12/01"
# for every block of rows in the sheet that starts a new company's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025/11/30") {
        # Force a text assignment (avoid Excel auto-parsing the string as a
        # date serial number), then clear the formatting change so the
        # cell's style stays exactly as it was (no explicit style index).
        $cell.NumberFormat = "@"
        $cell.Value = "2025/12/01"
        $cell.ClearFormats()
    }
}
